# Applies the updated coin price/volume figures (and the row-50/51
# SynthetixNetwork/Aptos refresh) captured in the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as literal text (not auto-converted
# to a number/date by Excel), matching the inlineStr cells in the source file.
function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

Set-TextValue 'D2' '29.204.34'
Set-TextValue 'E2' '  -1.15%  '
Set-TextValue 'D3' '1.859.97'
Set-TextValue 'E3' '  -0.89%  '
Set-TextValue 'D4' '0.9996'
Set-TextValue 'E4' '  -0.12%  '
Set-TextValue 'D5' '0.7145'
Set-TextValue 'E5' '  -1.05%  '
Set-TextValue 'D6' '240.46'
Set-TextValue 'E6' '  +0.22%  '
Set-TextValue 'D7' '0.9995'
Set-TextValue 'E7' '  -0.15%  '
Set-TextValue 'D8' '0.07753'
Set-TextValue 'E8' '  -1.24%  '
Set-TextValue 'D9' '0.3080'
Set-TextValue 'E9' '  -0.21%  '
Set-TextValue 'D10' '25.05'
Set-TextValue 'E10' '  -1.10%  '
Set-TextValue 'D11' '0.08254'
Set-TextValue 'E11' '  +0.13%  '
Set-TextValue 'D12' '1.855.10'
Set-TextValue 'E12' '  -1.65%  '
Set-TextValue 'D13' '5.234'
Set-TextValue 'E13' '  -0.36%  '
Set-TextValue 'D14' '0.7157'
Set-TextValue 'E14' '  -1.47%  '
Set-TextValue 'D15' '90.23'
Set-TextValue 'E15' '  +0.09%  '
Set-TextValue 'D16' '29.181.06'
Set-TextValue 'E16' '  -1.46%  '
Set-TextValue 'D17' '5.871'
Set-TextValue 'E17' '  +0.39%  '
Set-TextValue 'D18' '244.55'
Set-TextValue 'E18' '  +0.52%  '
Set-TextValue 'D19' '0.000007808'
Set-TextValue 'E19' '  -0.65%  '
Set-TextValue 'D20' '13.16'
Set-TextValue 'E20' '  -1.47%  '
Set-TextValue 'D21' '2.113.77'
Set-TextValue 'E21' '  -1.59%  '
Set-TextValue 'D22' '0.9997'
Set-TextValue 'E22' '  -0.07%  '
Set-TextValue 'D23' '7.955'
Set-TextValue 'E23' '  +1.90%  '
Set-TextValue 'D24' '0.9996'
Set-TextValue 'E24' '  -0.14%  '
Set-TextValue 'D25' '0.1588'
Set-TextValue 'E25' '  +0.58%  '
Set-TextValue 'D26' '162.66'
Set-TextValue 'E26' '  +0.02%  '
Set-TextValue 'D27' '8.927'
Set-TextValue 'E27' '  -0.78%  '
Set-TextValue 'E28' '  -0.44%  '
Set-TextValue 'D29' '1.495'
Set-TextValue 'E29' '  +0.45%  '
Set-TextValue 'E30' '  -2.96%  '
Set-TextValue 'D31' '4.386'
Set-TextValue 'E31' '  +0.89%  '
Set-TextValue 'D32' '4.151'
Set-TextValue 'E32' '  +1.78%  '
Set-TextValue 'D33' '0.05190'
Set-TextValue 'E33' '  -1.25%  '
Set-TextValue 'E34' '  -1.95%  '
Set-TextValue 'D35' '1.173'
Set-TextValue 'E35' '  -2.23%  '
Set-TextValue 'D36' '0.7277'
Set-TextValue 'E36' '  +1.14%  '
Set-TextValue 'D37' '2.678'
Set-TextValue 'E37' '  +0.25%  '
Set-TextValue 'D38' '0.01854'
Set-TextValue 'E38' '  -0.95%  '
Set-TextValue 'D39' '2.686'
Set-TextValue 'E39' '  -1.12%  '
Set-TextValue 'D40' '1.158.75'
Set-TextValue 'E40' '  -2.02%  '
Set-TextValue 'D41' '0.9028'
Set-TextValue 'E41' '  -0.99%  '
Set-TextValue 'D42' '6.093'
Set-TextValue 'E42' '  +1.48%  '
Set-TextValue 'E43' '  +0.83%  '
Set-TextValue 'D44' '0.9990'
Set-TextValue 'E44' '  -0.20%  '
Set-TextValue 'D45' '101.68'
Set-TextValue 'E45' '  -1.38%  '
Set-TextValue 'D46' '2.006.66'
Set-TextValue 'E46' '  -1.67%  '
Set-TextValue 'D47' '0.5230'
Set-TextValue 'E47' '  -2.34%  '
Set-TextValue 'D48' '1.769'
Set-TextValue 'E48' '  -0.69%  '
Set-TextValue 'D49' '9.307'
Set-TextValue 'E49' '  +0.74%  '
Set-TextValue 'B50' 'SynthetixNetwork'
Set-TextValue 'C50' 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextValue 'D50' '2.870'
Set-TextValue 'E50' '  +0.33%  '
Set-TextValue 'B51' 'Aptos'
Set-TextValue 'C51' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D51' '7.055'
Set-TextValue 'E51' '  -0.31%  '
